# Fix the typo'd duplicate shared string "Korper merkmale" (missing umlaut)
# in cell A13 so that it matches the correctly-spelled "Körper Merkmale"
# used elsewhere (e.g. A4). Once no cell references the old misspelled
# string, Excel will drop it from the shared-strings table on save,
# shifting the indices of all later shared strings down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A13").Value = "Körper Merkmale"

# Update the active selection to A13, matching the recorded cursor move.
$ws.Range("A13").Select()
